# Adds the three missing "logout" follow-up commands to the users-DB
# example tables (Input column): the exercise now shows which user each
# logout call applies to ("logout Simo" / "logout pesho" / "logout Simo").
#
# There are 6 occurrences of the word "logout" in the whole document; the
# first 3 are prose (explanatory text), the last 3 are the three example
# lines inside the two "Input" table cells that must receive the new
# username suffix.

$d = $word.ActiveDocument

$insertions = @{
    4 = " Simo"
    5 = " pesho"
    6 = " Simo"
}

$searchStart = 0
$matchIndex = 0

while ($true) {
    $matchIndex = $matchIndex + 1

    $rng = $d.Content
    $rng.Start = $searchStart
    $found = $rng.Find.Execute("logout", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    if (-not $found) {
        break
    }

    $nextStart = $rng.End

    if ($insertions.ContainsKey($matchIndex)) {
        $rng.Collapse(0)
        $rng.InsertAfter($insertions[$matchIndex])
        $nextStart = $rng.End
    }

    $searchStart = $nextStart
}
